# Apply cryptos list refresh: update prices / 1h volume %, and reorder a few rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    # Force the cell to remain a text/string cell (matches the source data,
    # which stores every value - including numeric-looking prices - as text),
    # then restore the default "Normal" style so no stray formatting is left behind.
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell "D2" "44.286.70"
Set-TextCell "E2" "  +0.45%  "

# Row 3
Set-TextCell "D3" "2.242.29"
Set-TextCell "E3" "  -0.14%  "

# Row 4
Set-TextCell "D4" "1.01"
Set-TextCell "E4" "  +0.26%  "

# Row 5
Set-TextCell "D5" "307.40"
Set-TextCell "E5" "  -2.65%  "

# Row 6
Set-TextCell "D6" "94.48"
Set-TextCell "E6" "  -4.96%  "

# Row 7
Set-TextCell "D7" "0.571"
Set-TextCell "E7" "  -0.64%  "

# Row 8
Set-TextCell "E8" "  +0.27%  "

# Row 9
Set-TextCell "D9" "0.525"
Set-TextCell "E9" "  -1.47%  "

# Row 10
Set-TextCell "D10" "34.78"
Set-TextCell "E10" "  -3.88%  "

# Row 11
Set-TextCell "D11" "0.0812"
Set-TextCell "E11" "  -1.46%  "

# Row 12
Set-TextCell "D12" "7.19"
Set-TextCell "E12" "  -2.37%  "

# Row 13
Set-TextCell "E13" "  -0.01%  "

# Row 14
Set-TextCell "D14" "2.582.34"
Set-TextCell "E14" "  -0.11%  "

# Row 15
Set-TextCell "D15" "2.240.93"
Set-TextCell "E15" "  -0.15%  "

# Row 16
Set-TextCell "D16" "0.831"
Set-TextCell "E16" "  -1.48%  "

# Row 17
Set-TextCell "D17" "13.52"
Set-TextCell "E17" "  -3.30%  "

# Row 18
Set-TextCell "D18" "44.029.07"
Set-TextCell "E18" "  +0.19%  "

# Row 19
Set-TextCell "D19" "0.0₃0964"
Set-TextCell "E19" "  -1.50%  "

# Row 20
Set-TextCell "D20" "6.40"
Set-TextCell "E20" "  +1.21%  "

# Row 21
Set-TextCell "D21" "12.14"
Set-TextCell "E21" "  -7.85%  "

# Row 22
Set-TextCell "D22" "65.56"
Set-TextCell "E22" "  -0.09%  "

# Row 23
Set-TextCell "D23" "238.28"
Set-TextCell "E23" "  +0.69%  "

# Row 24
Set-TextCell "D24" "2.95"
Set-TextCell "E24" "  -1.08%  "

# Row 25
Set-TextCell "D25" "2.00"
Set-TextCell "E25" "  -0.83%  "

# Row 26
Set-TextCell "E26" "  +0.09%  "

# Row 27
Set-TextCell "B27" "InjectiveProtocol"
Set-TextCell "C27" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell "D27" "38.53"
Set-TextCell "E27" "  +5.43%  "

# Row 28
Set-TextCell "B28" "Cosmos"
Set-TextCell "C28" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell "D28" "9.93"
Set-TextCell "E28" "  -2.48%  "

# Row 29
Set-TextCell "B29" "Toncoin"
Set-TextCell "C29" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell "D29" "2.21"
Set-TextCell "E29" "  +3.83%  "

# Row 30
Set-TextCell "D30" "20.03"
Set-TextCell "E30" "  -0.27%  "

# Row 31
Set-TextCell "D31" "5.85"
Set-TextCell "E31" "  -2.26%  "

# Row 32
Set-TextCell "D32" "153.07"
Set-TextCell "E32" "  -1.66%  "

# Row 33
Set-TextCell "D33" "0.0795"
Set-TextCell "E33" "  -5.11%  "

# Row 34
Set-TextCell "E34" "  -2.02%  "

# Row 35
Set-TextCell "E35" "  -5.57%  "

# Row 36
Set-TextCell "E36" "  +1.91%  "

# Row 37
Set-TextCell "D37" "0.107"
Set-TextCell "E37" "  -1.06%  "

# Row 38
Set-TextCell "D38" "1.77"
Set-TextCell "E38" "  -7.42%  "

# Row 39
Set-TextCell "D39" "3.52"
Set-TextCell "E39" "  -0.32%  "

# Row 40
Set-TextCell "D40" "3.81"
Set-TextCell "E40" "  -4.62%  "

# Row 41
Set-TextCell "D41" "14.33"
Set-TextCell "E41" "  -8.03%  "

# Row 42
Set-TextCell "D42" "0.0300"
Set-TextCell "E42" "  -2.38%  "

# Row 43
Set-TextCell "E43" "  +0.28%  "

# Row 44
Set-TextCell "D44" "1.751.80"
Set-TextCell "E44" "  +2.87%  "

# Row 45
Set-TextCell "D45" "82.97"
Set-TextCell "E45" "  +0.06%  "

# Row 46
Set-TextCell "D46" "0.191"
Set-TextCell "E46" "  -2.01%  "

# Row 47
Set-TextCell "D47" "99.88"
Set-TextCell "E47" "  -1.98%  "

# Row 48
Set-TextCell "D48" "4.93"
Set-TextCell "E48" "  -4.92%  "

# Row 49
Set-TextCell "D49" "8.09"
Set-TextCell "E49" "  -0.58%  "

# Row 50
Set-TextCell "B50" "MultiversX"
Set-TextCell "C50" "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextCell "D50" "54.76"
Set-TextCell "E50" "  -2.86%  "

# Row 51
Set-TextCell "B51" "Stacks"
Set-TextCell "C51" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell "D51" "1.57"
Set-TextCell "E51" "  -2.06%  "
